$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates derived from the cryptos list refresh (price / 1h volume
# changes, plus a newly inserted "PaxosStandard" row that shifted the rows
# below it down by one, pushing "RenderToken" off the bottom of the table).
$updates = @(
    @{ Cell = 'D2'; Value = '26.959.73' },
    @{ Cell = 'E2'; Value = '  +1.96%  ' },
    @{ Cell = 'D3'; Value = '1.815.26' },
    @{ Cell = 'E3'; Value = '  +2.38%  ' },
    @{ Cell = 'E4'; Value = '  +0.30%  ' },
    @{ Cell = 'D5'; Value = '312.27' },
    @{ Cell = 'E5'; Value = '  +1.85%  ' },
    @{ Cell = 'E6'; Value = '  +0.29%  ' },
    @{ Cell = 'D7'; Value = '0.4300' },
    @{ Cell = 'E7'; Value = '  +0.01%  ' },
    @{ Cell = 'D8'; Value = '0.3673' },
    @{ Cell = 'E8'; Value = '  -0.02%  ' },
    @{ Cell = 'D9'; Value = '0.07262' },
    @{ Cell = 'E9'; Value = '  +0.38%  ' },
    @{ Cell = 'D10'; Value = '2.153.18' },
    @{ Cell = 'E10'; Value = '  +21.26%  ' },
    @{ Cell = 'D11'; Value = '0.8654' },
    @{ Cell = 'E11'; Value = '  +1.94%  ' },
    @{ Cell = 'D12'; Value = '21.33' },
    @{ Cell = 'E12'; Value = '  +4.97%  ' },
    @{ Cell = 'D13'; Value = '5.408' },
    @{ Cell = 'E13'; Value = '  +3.04%  ' },
    @{ Cell = 'D14'; Value = '6.607' },
    @{ Cell = 'E14'; Value = '  +2.56%  ' },
    @{ Cell = 'D15'; Value = '0.06941' },
    @{ Cell = 'E15'; Value = '  +1.21%  ' },
    @{ Cell = 'D16'; Value = '81.00' },
    @{ Cell = 'E16'; Value = '  +1.71%  ' },
    @{ Cell = 'D17'; Value = '1.006' },
    @{ Cell = 'E17'; Value = '  +0.00%  ' },
    @{ Cell = 'D18'; Value = '0.000008877' },
    @{ Cell = 'E18'; Value = '  +2.12%  ' },
    @{ Cell = 'D19'; Value = '1.006' },
    @{ Cell = 'E19'; Value = '  +0.30%  ' },
    @{ Cell = 'D20'; Value = '15.26' },
    @{ Cell = 'E20'; Value = '  +1.34%  ' },
    @{ Cell = 'D21'; Value = '26.993.64' },
    @{ Cell = 'E21'; Value = '  +2.09%  ' },
    @{ Cell = 'D22'; Value = '5.187' },
    @{ Cell = 'E22'; Value = '  +1.50%  ' },
    @{ Cell = 'D23'; Value = '11.01' },
    @{ Cell = 'E23'; Value = '  -2.50%  ' },
    @{ Cell = 'D24'; Value = '2.358.32' },
    @{ Cell = 'E24'; Value = '  +18.14%  ' },
    @{ Cell = 'D25'; Value = '153.94' },
    @{ Cell = 'E25'; Value = '  +1.02%  ' },
    @{ Cell = 'E26'; Value = '  +1.68%  ' },
    @{ Cell = 'D27'; Value = '18.35' },
    @{ Cell = 'E27'; Value = '  +1.03%  ' },
    @{ Cell = 'D28'; Value = '5.223' },
    @{ Cell = 'E28'; Value = '  +2.52%  ' },
    @{ Cell = 'D29'; Value = '1.900' },
    @{ Cell = 'E29'; Value = '  +10.53%  ' },
    @{ Cell = 'D30'; Value = '114.65' },
    @{ Cell = 'D31'; Value = '0.08947' },
    @{ Cell = 'E31'; Value = '  +0.17%  ' },
    @{ Cell = 'D32'; Value = '1.188' },
    @{ Cell = 'E32'; Value = '  +6.19%  ' },
    @{ Cell = 'D33'; Value = '0.7464' },
    @{ Cell = 'E33'; Value = '  +3.00%  ' },
    @{ Cell = 'E34'; Value = '  +1.88%  ' },
    @{ Cell = 'D35'; Value = '2.810' },
    @{ Cell = 'E35'; Value = '  +2.17%  ' },
    @{ Cell = 'E36'; Value = '  +0.29%  ' },
    @{ Cell = 'D37'; Value = '1.131' },
    @{ Cell = 'E37'; Value = '  +4.89%  ' },
    @{ Cell = 'D38'; Value = '0.05214' },
    @{ Cell = 'E38'; Value = '  +0.97%  ' },
    @{ Cell = 'D39'; Value = '0.01925' },
    @{ Cell = 'E39'; Value = '  +1.44%  ' },
    @{ Cell = 'D40'; Value = '0.5100' },
    @{ Cell = 'E40'; Value = '  +3.49%  ' },
    @{ Cell = 'D41'; Value = '0.1656' },
    @{ Cell = 'E41'; Value = '  +3.00%  ' },
    @{ Cell = 'D42'; Value = '2.736' },
    @{ Cell = 'E42'; Value = '  +7.96%  ' },
    @{ Cell = 'D43'; Value = '6.472' },
    @{ Cell = 'E43'; Value = '  +4.13%  ' },
    @{ Cell = 'D44'; Value = '8.311' },
    @{ Cell = 'E44'; Value = '  +3.32%  ' },
    @{ Cell = 'B45'; Value = 'PaxosStandard' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax' },
    @{ Cell = 'D45'; Value = '1.006' },
    @{ Cell = 'E45'; Value = '  +0.30%  ' },
    @{ Cell = 'B46'; Value = 'Quant' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Cell = 'D46'; Value = '106.77' },
    @{ Cell = 'E46'; Value = '  +1.77%  ' },
    @{ Cell = 'B47'; Value = 'EnergySwap' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D47'; Value = '10.40' },
    @{ Cell = 'E47'; Value = '  +2.26%  ' },
    @{ Cell = 'B48'; Value = 'PaxDollar' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Cell = 'D48'; Value = '1.005' },
    @{ Cell = 'E48'; Value = '  +0.34%  ' },
    @{ Cell = 'B49'; Value = 'Decentraland' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' },
    @{ Cell = 'D49'; Value = '0.4589' },
    @{ Cell = 'E49'; Value = '  +2.22%  ' },
    @{ Cell = 'B50'; Value = 'NEARProtocol' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = 'D50'; Value = '1.644' },
    @{ Cell = 'E50'; Value = '  +3.91%  ' },
    @{ Cell = 'B51'; Value = 'Cronos' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'D51'; Value = '0.06214' },
    @{ Cell = 'E51'; Value = '  +0.25%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "81.00", "1.006")
    # are preserved exactly instead of being coerced into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Restore the default "Normal" style so no stray number-format style
    # gets attached to the cell (matches original inline-string cells,
    # which carry no explicit style).
    $cell.Style = "Normal"
}
